$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a trailing space to each question in column A (rows 2-10) so they
# match the "What is your birth month? " style already used on row 3.
# Row 1 (header "Questions"/"DataType") and column B ("String") are untouched.
$ws.Range("A2").Value = "What is your name? "
$ws.Range("A4").Value = "What is your favorite animal? "
$ws.Range("A5").Value = "What is your favorite plant? "
$ws.Range("A6").Value = "What is your hobby? "
$ws.Range("A7").Value = "What is your favorite food? "
$ws.Range("A8").Value = "What is your favorite colour? "
$ws.Range("A9").Value = "What is your favorite music genre? "
$ws.Range("A10").Value = "What is your favorite part of the day? "

# Move the active selection to A4 (matches the saved workbook's cursor state).
$ws.Range("A4").Select() | Out-Null
